# Update Sheet1 with all matches from 2017/18 season.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update points_won (D), points_lost (E) and place (F) columns for the
# rows belonging to the 2017/2018 season.
$ws.Range("D2").Value = 37
$ws.Range("E2").Value = 19
$ws.Range("F2").Value = 4

$ws.Range("D10").Value = 40
$ws.Range("E10").Value = 4

$ws.Range("D18").Value = 38
$ws.Range("E18").Value = 6

$ws.Range("D26").Value = 17
$ws.Range("E26").Value = 27

$ws.Range("D34").Value = 24
$ws.Range("E34").Value = 16

$ws.Range("D42").Value = 20
$ws.Range("E42").Value = 16
$ws.Range("F42").Value = 4

# Update the saved view state (scroll position and selection) to match
# where the editor left off after making the updates.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 21
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F34").Select()
